# Applies the OOXML diff to the document:
#  1. Wrap "Description" table-header text in proofErr spellStart/spellEnd.
#  2. Split "IFileSystemComponent" out of its sentence into its own run,
#     wrapped in proofErr spellStart/spellEnd (two occurrences).
#  3. Delete the "Information Expert" ... "Controller & Indirection" section
#     (six paragraphs), collapsing down to a single empty paragraph that
#     keeps the trailing _GoBack bookmark.

$d = $word.ActiveDocument

function Set-ParagraphXml($paragraph, [string]$innerXml) {
    $range = $paragraph.Range.Duplicate()
    $pkg = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>" + $innerXml + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    $range.InsertXML($pkg)
}

function Get-ParaText($paragraph) {
    return $paragraph.Range.Text.TrimEnd([char]13, [char]7)
}

# --- 1. "Description" table header cell: wrap run in proofErr markers ---
$descPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ((Get-ParaText $cand) -eq "Description") {
        $descPara = $cand
        break
    }
}
if ($descPara -ne $null) {
    $descXml = "<w:p><w:proofErr w:type='spellStart'/><w:r><w:t>Description</w:t></w:r><w:proofErr w:type='spellEnd'/></w:p>"
    Set-ParagraphXml $descPara $descXml
}

# --- 2a. "We will do this ... IFileSystemComponent ..." paragraph ---
$quote = [char]0x201C
$rquote = [char]0x201D
$needle1 = "We will do this by having both the document and the folder class inherit from the interface " + $quote + "IFileSystemComponent" + $rquote + ". "
$p1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ((Get-ParaText $cand) -eq $needle1) {
        $p1 = $cand
        break
    }
}
if ($p1 -ne $null) {
    $xml1 = "<w:p><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" `
        + "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>We will do this by having both the document and the folder class inherit from the interface " + $quote + "</w:t></w:r>" `
        + "<w:proofErr w:type='spellStart'/>" `
        + "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>IFileSystemComponent</w:t></w:r>" `
        + "<w:proofErr w:type='spellEnd'/>" `
        + "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>" + $rquote + ". </w:t></w:r>" `
        + "</w:p>"
    Set-ParagraphXml $p1 $xml1
}

# --- 2b. "A folder will then contain a list of IFileSystemComponent's ..." paragraph ---
$apos = [char]0x2019
$needle2 = "A folder will then contain a list of IFileSystemComponent" + $apos + "s which then obviously can be both documents and other folders."
$p2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ((Get-ParaText $cand) -eq $needle2) {
        $p2 = $cand
        break
    }
}
if ($p2 -ne $null) {
    $xml2 = "<w:p><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" `
        + "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>A folder will then contain a list of </w:t></w:r>" `
        + "<w:proofErr w:type='spellStart'/>" `
        + "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>IFileSystemComponent" + $apos + "s</w:t></w:r>" `
        + "<w:proofErr w:type='spellEnd'/>" `
        + "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> which then obviously can be both documents and other folders.</w:t></w:r>" `
        + "</w:p>"
    Set-ParagraphXml $p2 $xml2
}

# --- 3. Remove "Information Expert" .. "Controller & Indirection" section ---
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ((Get-ParaText $cand) -eq "Information Expert") {
        $startPara = $cand
        break
    }
}
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ((Get-ParaText $cand) -eq "As a controller our controller takes the users input and then calls relevant functions in other modules like our information expert, the Storage.") {
        $endPara = $cand
    }
}

if (($startPara -ne $null) -and ($endPara -ne $null)) {
    $bigRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $xml3 = "<w:p><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>"
    $pkg3 = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>" + $xml3 + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    $bigRange.InsertXML($pkg3)
}

Write-Output "edit complete"
